$d = $word.ActiveDocument

# 1. Fix "opoosing" -> "opposing" (also removes surrounding proofErr spell tags naturally
#    since Find/Replace rewrites the run text)
$d.Content.Find.Execute("opoosing", $false, $false, $false, $false, $false, $true, 1, $false, "opposing", 2)

Write-Host "Done"
